$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08474199241993637
$ws.Range("H2").Value = 1.042420173395213
$ws.Range("I2").Value = 23.59358390711208
$ws.Range("G3").Value = 0.1232941647648896
$ws.Range("H3").Value = 5.96313129393865
$ws.Range("G4").Value = -0.02606606164238651
$ws.Range("H4").Value = -68.98275899777352
$ws.Range("G5").Value = -0.01036525652738265
$ws.Range("H5").Value = 85.57409181188447
$ws.Range("G6").Value = 0.03386128906426435
$ws.Range("H6").Value = -4.738625622107443
$ws.Range("G7").Value = 0.04884726572672313
$ws.Range("H7").Value = 140.6850410674777
$ws.Range("G8").Value = -0.1498574066992334
$ws.Range("H8").Value = -6.208108098544041
$ws.Range("G9").Value = -0.1342423540895483
$ws.Range("H9").Value = 1.931999336120997
$ws.Range("G10").Value = -0.06389265820825843
$ws.Range("H10").Value = 40.85498497729278
$ws.Range("G11").Value = -0.07570387819327916
$ws.Range("H11").Value = -13.77051082700729
$ws.Range("G12").Value = -0.355145767195666
$ws.Range("H12").Value = 14.46317147635376
$ws.Range("G13").Value = -0.38328291969247
$ws.Range("H13").Value = 14.61718819232252
$ws.Range("G14").Value = -0.04731702611729585
$ws.Range("H14").Value = 6.738194534252047
$ws.Range("G15").Value = -0.03398353532849006
$ws.Range("H15").Value = 58.90319917783064
$ws.Range("G16").Value = 0.1176299041054609
$ws.Range("H16").Value = -19.44408948185969
$ws.Range("G17").Value = 0.1696191756508549
$ws.Range("H17").Value = 38.32745155380473
$ws.Range("G18").Value = 0.1257455380488698
$ws.Range("H18").Value = -9.180919715594834
$ws.Range("G19").Value = 0.1554091097657176
$ws.Range("H19").Value = 62.89111433033775
$ws.Range("G20").Value = 0.02851946484408382
$ws.Range("H20").Value = 11.25617082861953
$ws.Range("G21").Value = 0.0571340621353415
$ws.Range("H21").Value = -23.99061209759593
$ws.Range("G24").Value = 0.08777057322145604
$ws.Range("H24").Value = -12.61088635272248
$ws.Range("G25").Value = 0.1627668706529699
$ws.Range("H25").Value = 7.403888662569121
$ws.Range("G26").Value = 0.07314343790096881
$ws.Range("H26").Value = -7.559299080282517
$ws.Range("G27").Value = 0.09029633928040119
$ws.Range("H27").Value = -9.570386529117695
$ws.Range("G28").Value = -0.2418225063297041
$ws.Range("H28").Value = -13.46820175375412
$ws.Range("G29").Value = -0.2096819259979749
$ws.Range("H29").Value = -2.13108031294146
$ws.Range("G30").Value = 0.04901987756793426
$ws.Range("H30").Value = 11.07648529888407
$ws.Range("G31").Value = 0.03159286549917024
$ws.Range("H31").Value = 19.97391601840917
$ws.Range("G32").Value = 0.1000229799695457
$ws.Range("H32").Value = 5.331133195689698
$ws.Range("G33").Value = 0.1398811186130419
$ws.Range("H33").Value = 34.54650190954391
$ws.Range("G34").Value = 0.04629473990802215
$ws.Range("H34").Value = -0.2863622768830413
$ws.Range("G35").Value = 0.02024525192158527
$ws.Range("H35").Value = 167.2101672157798
$ws.Range("G36").Value = 0.06121738559377108
$ws.Range("H36").Value = 6.029383464062114
$ws.Range("G37").Value = 0.09505657619121473
$ws.Range("H37").Value = 35.16531692067698
$ws.Range("G38").Value = 0.02786707858175142
$ws.Range("H38").Value = -46.80475746705679
$ws.Range("G39").Value = 0.0415705448348859
$ws.Range("H39").Value = 100.4623844331303
$ws.Range("G40").Value = 0.005199291201779043
$ws.Range("H40").Value = 161.2330422624623
$ws.Range("G41").Value = 0.03048685014122112
$ws.Range("H41").Value = -13.7714220347054
$ws.Range("G42").Value = 0.1439807322537644
$ws.Range("H42").Value = 7.700579553901894
$ws.Range("G43").Value = 0.1485117435490588
$ws.Range("H43").Value = -0.3118555327506884
$ws.Range("G44").Value = -0.004225429036736303
$ws.Range("H44").Value = 50.35612601282124
$ws.Range("G45").Value = -0.001142573481500866
$ws.Range("H45").Value = 89.59392055111434
$ws.Range("G46").Value = 0.007570830599660347
$ws.Range("H46").Value = 329.9035849951198
$ws.Range("G47").Value = -0.02170929581755342
$ws.Range("H47").Value = -133.9746731680467
$ws.Range("G48").Value = 0.04951042073931265
$ws.Range("H48").Value = -1.515136807752838
$ws.Range("G49").Value = 0.07308149737294023
$ws.Range("H49").Value = 10.62095923674684
$ws.Range("G50").Value = 0.1699639866851714
$ws.Range("H50").Value = 5.399002608764649
$ws.Range("G51").Value = 0.1456389746321311
$ws.Range("H51").Value = -14.88822772611109
$ws.Range("G52").Value = -0.1664366200929072
$ws.Range("H52").Value = -3.743427651143814
$ws.Range("G53").Value = -0.1287061054412911
$ws.Range("H53").Value = -2.10353861198138
$ws.Range("G54").Value = 0.1324501758550223
$ws.Range("H54").Value = 41.32179995338848
$ws.Range("G55").Value = 0.1250001293916364
$ws.Range("H55").Value = 10.54059431180212
$ws.Range("G56").Value = -0.02769906760815561
$ws.Range("H56").Value = -279.3991852999104
$ws.Range("G57").Value = -0.02168662540772742
$ws.Range("H57").Value = 5.163535363577587
$ws.Range("G58").Value = 0.06698151312512451
$ws.Range("H58").Value = 18.78687968529585
$ws.Range("G59").Value = 0.06585844772891644
$ws.Range("H59").Value = -8.300590902478756
$ws.Range("G60").Value = 0.05890242968896087
$ws.Range("H60").Value = -15.82234990371945
$ws.Range("G61").Value = 0.07271139302061902
$ws.Range("H61").Value = 52.9903447276155
$ws.Range("G62").Value = 0.07626640707222895
$ws.Range("H62").Value = 4.531338705753103
$ws.Range("G63").Value = 0.06804087353861099
$ws.Range("H63").Value = 4.048366905193829
$ws.Range("G64").Value = -0.0317303215415714
$ws.Range("H64").Value = 23.38540944257976
$ws.Range("G65").Value = 0.009559373620914576
$ws.Range("H65").Value = 119.3808901814884
$ws.Range("G66").Value = 0.02716256529430541
$ws.Range("H66").Value = 43.45938826848288
$ws.Range("G67").Value = 0.01961745882999849
$ws.Range("H67").Value = -24.9807095175409
$ws.Range("G68").Value = 0.003938328107505674
$ws.Range("H68").Value = 590.9912760576408
$ws.Range("G69").Value = 0.01171122498132751
$ws.Range("H69").Value = 190.6981623774254
$ws.Range("G70").Value = -0.03948035243566413
$ws.Range("H70").Value = -43.83654349643278
$ws.Range("G71").Value = -0.03306707538665055
$ws.Range("H71").Value = 39.97839104885473
$ws.Range("G72").Value = -0.1293312793781893
$ws.Range("H72").Value = 12.81085123638627
$ws.Range("G73").Value = -0.1461414291459345
$ws.Range("H73").Value = -0.933486203242633
$ws.Range("G74").Value = 0.146236662701106
$ws.Range("H74").Value = 16.01378818809507
$ws.Range("G75").Value = 0.146585111236492
$ws.Range("H75").Value = 8.443417234805935
$ws.Range("G76").Value = -0.04758209433378748
$ws.Range("H76").Value = -38.16763537144495
$ws.Range("G77").Value = -0.03536116151610885
$ws.Range("H77").Value = 23.43915925703632
$ws.Range("G78").Value = 0.07833183570200634
$ws.Range("H78").Value = -15.01515841924144
$ws.Range("G79").Value = 0.1038733755289636
$ws.Range("H79").Value = 7.635863034048246
$ws.Range("G80").Value = -0.1925323710809149
$ws.Range("H80").Value = -18.54502961563816
$ws.Range("G81").Value = -0.1851666689320838
$ws.Range("H81").Value = 14.45230396767276
$ws.Range("G82").Value = 0.1630484570105964
$ws.Range("H82").Value = 17.51555291012132
$ws.Range("G83").Value = 0.1947080810003642
$ws.Range("H83").Value = 18.277584770784
$ws.Range("G84").Value = 0.007823077013249097
$ws.Range("H84").Value = -44.101024598406
$ws.Range("G85").Value = 0.03253500298388964
$ws.Range("H85").Value = 43.71102513562325
